$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.7596665024757385
$ws.Range("B1").Value = 1.110942602157593
$ws.Range("C1").Value = 3.503080368041992
$ws.Range("D1").Value = 3.807764291763306
$ws.Range("E1").Value = 1.977208971977234
